# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F8").Value = 11514
    $ws.Range("F14").Value = 803
    $ws.Range("F16").Value = 13083
    $ws.Range("F24").Value = 120
}
